function Set-TextValue {
    param($ws, $address, $text)
    $helper = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $ws.Range($address).PasteSpecial(-4163)
    $helper.Clear()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '71.999.22'
$ws.Range('E2').Value = '  +3.74%  '
Set-TextValue $ws 'D3' '3.598.73'
$ws.Range('E3').Value = '  +6.20%  '
Set-TextValue $ws 'D4' '1.00'
$ws.Range('E4').Value = '  +0.12%  '
Set-TextValue $ws 'D5' '594.07'
$ws.Range('E5').Value = '  +1.16%  '
Set-TextValue $ws 'D6' '184.02'
$ws.Range('E6').Value = '  +2.38%  '
Set-TextValue $ws 'D7' '3.591.56'
$ws.Range('E7').Value = '  +6.07%  '
Set-TextValue $ws 'D8' '0.608'
$ws.Range('E8').Value = '  +1.85%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  +6.41%  '
Set-TextValue $ws 'D11' '0.609'
$ws.Range('E11').Value = '  +2.84%  '
Set-TextValue $ws 'D12' '50.24'
$ws.Range('E12').Value = '  +3.39%  '
Set-TextValue $ws 'D13' '0.0000290'
$ws.Range('E13').Value = '  +2.71%  '
Set-TextValue $ws 'D14' '700.89'
$ws.Range('E14').Value = '  +3.29%  '
Set-TextValue $ws 'D15' '4.167.52'
$ws.Range('E15').Value = '  +5.97%  '
Set-TextValue $ws 'D16' '8.91'
$ws.Range('E16').Value = '  +3.29%  '
Set-TextValue $ws 'D17' '72.010.53'
$ws.Range('E17').Value = '  +3.65%  '
Set-TextValue $ws 'D18' '3.563.78'
$ws.Range('E18').Value = '  +5.05%  '
$ws.Range('E19').Value = '  +1.32%  '
Set-TextValue $ws 'D20' '18.31'
$ws.Range('E20').Value = '  +3.50%  '
Set-TextValue $ws 'D21' '11.77'
$ws.Range('E21').Value = '  +4.23%  '
Set-TextValue $ws 'D22' '0.930'
$ws.Range('E22').Value = '  +2.74%  '
Set-TextValue $ws 'D23' '5.74'
$ws.Range('E23').Value = '  +6.06%  '
Set-TextValue $ws 'D24' '17.75'
$ws.Range('E24').Value = '  +3.16%  '
Set-TextValue $ws 'D25' '104.65'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('E26').Value = '  +2.17%  '
Set-TextValue $ws 'D27' '2.82'
$ws.Range('E27').Value = '  +3.16%  '
$ws.Range('E28').Value = '  +2.84%  '
Set-TextValue $ws 'D29' '35.37'
$ws.Range('E29').Value = '  +3.48%  '
$ws.Range('E30').Value = '  +3.77%  '
$ws.Range('E31').Value = '  +6.52%  '
$ws.Range('E32').Value = '  +16.53%  '
Set-TextValue $ws 'D33' '589.17'
$ws.Range('E33').Value = '  +5.61%  '
Set-TextValue $ws 'D34' '11.35'
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('E35').Value = '  +0.61%  '
Set-TextValue $ws 'D36' '59.82'
$ws.Range('E36').Value = '  +2.14%  '
$ws.Range('E37').Value = '  -0.05%  '
Set-TextValue $ws 'D38' '3.677.37'
$ws.Range('E38').Value = '  -0.38%  '
Set-TextValue $ws 'D39' '0.146'
$ws.Range('E39').Value = '  +4.85%  '
Set-TextValue $ws 'D40' '36.30'
$ws.Range('E40').Value = '  +2.44%  '
Set-TextValue $ws 'D41' '0.0₃0785'
$ws.Range('E41').Value = '  +12.27%  '
Set-TextValue $ws 'D42' '3.48'
$ws.Range('E42').Value = '  +6.46%  '
$ws.Range('E43').Value = '  +4.32%  '
Set-TextValue $ws 'D44' '0.0441'
$ws.Range('E44').Value = '  +4.18%  '
Set-TextValue $ws 'D45' '0.348'
$ws.Range('E45').Value = '  +2.38%  '
$ws.Range('E46').Value = '  +1.93%  '
Set-TextValue $ws 'D47' '2.77'
$ws.Range('E47').Value = '  +3.52%  '
$ws.Range('E48').Value = '  +4.62%  '
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('E50').Value = '  -0.31%  '
Set-TextValue $ws 'D51' '134.44'
$ws.Range('E51').Value = '  +0.93%  '
